$wb = $excel.ActiveWorkbook

# --- Italy sheet: copy of Norway, placed right after Norway ---
$norway = $wb.Worksheets.Item("Norway")
$norway.Copy($null, $norway)
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3443/T1971"

# --- Spain sheet: copy of Italy, placed right after Italy ---
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3442/T2129"

# Selections to match the target state
$spain.Activate() | Out-Null
$spain.Range("A6").Select() | Out-Null

$italy.Activate() | Out-Null
$italy.Range("B4").Select() | Out-Null
